$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows (models "null=median" and "null=mean") right after
#    the existing "RF, null=-1" row, pushing the "RF, null managed with RF"
#    row down from row 3 to row 5.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the two new rows.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "RF, null=median, no class balancing"
$ws.Range("C3").Value = 0.838847926
$ws.Range("D3").Value = 0.00432671898687496

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "RF, null=mean, no class balancing"
$ws.Range("C4").Value = 0.840236824150493
$ws.Range("D4").Value = 0.00455476898471562
$ws.Range("E4").Value = 0.841631

# The old row 3 ("RF, null managed with RF...") is now row 5 - its "model"
# index bumps from 2 to 3.
$ws.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 3. Re-style everything: Calibri / Arial Unicode MS -> DejaVu Sans, and the
#    two numeric columns (internal AUC avg / std) get a fixed 6-decimal
#    number format. Each distinct look is built once on a scratch cell and
#    then format-painted (Copy + PasteSpecial formats) onto its destination
#    range so we don't leave stray intermediate styles behind.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

# -- style 1: bold header font, still General format -----------------------
$ws.Range("Z1").Font.Name = "DejaVu Sans"
$ws.Range("Z1").Font.Size = 11
$ws.Range("Z1").Font.Bold = $true
$ws.Range("Z1").Copy()
$ws.Range("A1:E1").PasteSpecial($xlPasteFormats)

# -- style 2: plain body font, General format -------------------------------
$ws.Range("Z2").Font.Name = "DejaVu Sans"
$ws.Range("Z2").Font.Size = 11
$ws.Range("Z2").Copy()
$ws.Range("A2:B2").PasteSpecial($xlPasteFormats)
$ws.Range("A3:B4").PasteSpecial($xlPasteFormats)
$ws.Range("A5:B5").PasteSpecial($xlPasteFormats)

# -- style 3: numeric, 6dp, vertically centered ------------------------------
$ws.Range("Z3").Font.Name = "DejaVu Sans"
$ws.Range("Z3").Font.Size = 11
$ws.Range("Z3").NumberFormat = "0.000000"
$ws.Range("Z3").VerticalAlignment = -4108
$ws.Range("Z3").Copy()
$ws.Range("C2:D2").PasteSpecial($xlPasteFormats)

# -- style 4: numeric, 6dp, default (bottom) alignment -----------------------
$ws.Range("Z4").Font.Name = "DejaVu Sans"
$ws.Range("Z4").Font.Size = 11
$ws.Range("Z4").NumberFormat = "0.000000"
$ws.Range("Z4").Copy()
$ws.Range("E2").PasteSpecial($xlPasteFormats)
$ws.Range("E3").PasteSpecial($xlPasteFormats)
$ws.Range("D5").PasteSpecial($xlPasteFormats)
$ws.Range("E5").PasteSpecial($xlPasteFormats)

# -- style 5: numeric, 6dp, wrapped text -------------------------------------
$ws.Range("Z5").Font.Name = "DejaVu Sans"
$ws.Range("Z5").Font.Size = 11
$ws.Range("Z5").NumberFormat = "0.000000"
$ws.Range("Z5").WrapText = $true
$ws.Range("Z5").Copy()
$ws.Range("C3:D4").PasteSpecial($xlPasteFormats)
$ws.Range("C5").PasteSpecial($xlPasteFormats)

# -- style 6: default font, General format, wrapped text --------------------
$ws.Range("Z6").WrapText = $true
$ws.Range("Z6").Copy()
$ws.Range("E4").PasteSpecial($xlPasteFormats)

# -- drop the scratch cells used to build the templates ----------------------
$ws.Range("Z1:Z6").Clear()

# ---------------------------------------------------------------------------
# 4. Column widths / row heights.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 44.665
$ws.Columns.Item(3).ColumnWidth = 17.8317
$ws.Columns.Item(4).ColumnWidth = 15.0033
$ws.Columns.Item(5).ColumnWidth = 12.165

$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 14.95
$ws.Rows.Item(4).RowHeight = 14.9
$ws.Rows.Item(5).RowHeight = 14.95

# ---------------------------------------------------------------------------
# 5. Selection cursor matches the authored file.
# ---------------------------------------------------------------------------
$ws.Range("D8").Select()
